$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.023937109182996
$ws.Cells.Item(2, 4).Value = 1.028192664279999
$ws.Cells.Item(2, 5).Value = 1.027561751191078
$ws.Cells.Item(2, 6).Value = 1.034328969650337
$ws.Cells.Item(2, 9).Value = 1.031056677312877
$ws.Cells.Item(2, 10).Value = 1.029114530096275
$ws.Cells.Item(2, 11).Value = 1.031010195307241
$ws.Cells.Item(2, 12).Value = 1.030381119595177
$ws.Cells.Item(2, 13).Value = 1.037128757269461
$ws.Cells.Item(2, 14).Value = 1.030575990770102

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.025020265068463
$ws.Cells.Item(3, 4).Value = 1.028978776059872
$ws.Cells.Item(3, 5).Value = 1.028592414181948
$ws.Cells.Item(3, 6).Value = 1.035614902369841
$ws.Cells.Item(3, 9).Value = 1.031294541125095
$ws.Cells.Item(3, 10).Value = 1.029835686314637
$ws.Cells.Item(3, 11).Value = 1.031604249765691
$ws.Cells.Item(3, 12).Value = 1.031218931341415
$ws.Cells.Item(3, 13).Value = 1.038222585196461
$ws.Cells.Item(3, 14).Value = 1.031298171113008

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.025720946000182
$ws.Cells.Item(4, 4).Value = 1.029486937832025
$ws.Cells.Item(4, 5).Value = 1.029259498675032
$ws.Cells.Item(4, 6).Value = 1.036446773265837
$ws.Cells.Item(4, 9).Value = 1.031446536856227
$ws.Cells.Item(4, 10).Value = 1.030301601551921
$ws.Cells.Item(4, 11).Value = 1.031987485753661
$ws.Cells.Item(4, 12).Value = 1.031760629780278
$ws.Cells.Item(4, 13).Value = 1.03892960864476
$ws.Cells.Item(4, 14).Value = 1.03176474800337

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.026015466899278
$ws.Cells.Item(5, 4).Value = 1.029700447958536
$ws.Cells.Item(5, 5).Value = 1.029539983650587
$ws.Cells.Item(5, 6).Value = 1.036796442503283
$ws.Cells.Item(5, 9).Value = 1.031509976635423
$ws.Cells.Item(5, 10).Value = 1.030497300024079
$ws.Cells.Item(5, 11).Value = 1.032148321233604
$ws.Cells.Item(5, 12).Value = 1.031988259174042
$ws.Cells.Item(5, 13).Value = 1.039226661693735
$ws.Cells.Item(5, 14).Value = 1.031960724389806

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.026064915620821
$ws.Cells.Item(6, 4).Value = 1.029736290143208
$ws.Cells.Item(6, 5).Value = 1.029587080824179
$ws.Cells.Item(6, 6).Value = 1.036855150670019
$ws.Cells.Item(6, 9).Value = 1.031520601523029
$ws.Cells.Item(6, 10).Value = 1.030530148564763
$ws.Cells.Item(6, 11).Value = 1.032175309980485
$ws.Cells.Item(6, 12).Value = 1.032026473259847
$ws.Cells.Item(6, 13).Value = 1.039276527717913
$ws.Cells.Item(6, 14).Value = 1.031993619579185

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.025724881579952
$ws.Cells.Item(7, 4).Value = 1.029489791240775
$ws.Cells.Item(7, 5).Value = 1.029263246360946
$ws.Cells.Item(7, 6).Value = 1.036451445754685
$ws.Cells.Item(7, 9).Value = 1.03144738634617
$ws.Cells.Item(7, 10).Value = 1.030304217161084
$ws.Cells.Item(7, 11).Value = 1.031989635933315
$ws.Cells.Item(7, 12).Value = 1.031763671770366
$ws.Cells.Item(7, 13).Value = 1.038933578586068
$ws.Cells.Item(7, 14).Value = 1.031767367326998

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.024303207379277
$ws.Cells.Item(8, 4).Value = 1.028458438932744
$ws.Cells.Item(8, 5).Value = 1.027910032229829
$ws.Cells.Item(8, 6).Value = 1.034763601962217
$ws.Cells.Item(8, 9).Value = 1.031137461627849
$ws.Cells.Item(8, 10).Value = 1.02935839771674
$ws.Cells.Item(8, 11).Value = 1.031211198384715
$ws.Cells.Item(8, 12).Value = 1.030664349540807
$ws.Cells.Item(8, 13).Value = 1.037498579033087
$ws.Cells.Item(8, 14).Value = 1.030820204710557

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.021796512903834
$ws.Cells.Item(9, 4).Value = 1.026637201789162
$ws.Cells.Item(9, 5).Value = 1.025526831290987
$ws.Cells.Item(9, 6).Value = 1.031787686941411
$ws.Cells.Item(9, 9).Value = 1.030576645717869
$ws.Cells.Item(9, 10).Value = 1.027686205118285
$ws.Cells.Item(9, 11).Value = 1.029830627802637
$ws.Cells.Item(9, 12).Value = 1.028723956717258
$ws.Cells.Item(9, 13).Value = 1.034964067334035
$ws.Cells.Item(9, 14).Value = 1.029145637406809

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.020124297591927
$ws.Cells.Item(10, 4).Value = 1.025420451461497
$ws.Cells.Item(10, 5).Value = 1.023938900702906
$ws.Cells.Item(10, 6).Value = 1.029802469068214
$ws.Cells.Item(10, 9).Value = 1.030192891914797
$ws.Cells.Item(10, 10).Value = 1.02656766233932
$ws.Cells.Item(10, 11).Value = 1.028904280316654
$ws.Cells.Item(10, 12).Value = 1.027428152982919
$ws.Cells.Item(10, 13).Value = 1.033270357818118
$ws.Cells.Item(10, 14).Value = 1.028025506168799

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.019399937720371
$ws.Cells.Item(11, 4).Value = 1.02489297024439
$ws.Cells.Item(11, 5).Value = 1.023251508391517
$ws.Cells.Item(11, 6).Value = 1.028942518092808
$ws.Cells.Item(11, 9).Value = 1.030024378942446
$ws.Cells.Item(11, 10).Value = 1.026082425785628
$ws.Cells.Item(11, 11).Value = 1.028501743768798
$ws.Cells.Item(11, 12).Value = 1.026866525184603
$ws.Cells.Item(11, 13).Value = 1.032535983615893
$ws.Cells.Item(11, 14).Value = 1.027539580523543

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.019130834828998
$ws.Cells.Item(12, 4).Value = 1.024696946916111
$ws.Cells.Item(12, 5).Value = 1.02299620841103
$ws.Cells.Item(12, 6).Value = 1.028623040673978
$ws.Cells.Item(12, 9).Value = 1.029961433252173
$ws.Cells.Item(12, 10).Value = 1.02590205144498
$ws.Cells.Item(12, 11).Value = 1.028352009808723
$ws.Cells.Item(12, 12).Value = 1.02665783037202
$ws.Cells.Item(12, 13).Value = 1.032263054051764
$ws.Cells.Item(12, 14).Value = 1.027358950030643

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.019188560298009
$ws.Cells.Item(13, 4).Value = 1.024738998810412
$ws.Cells.Item(13, 5).Value = 1.023050969862378
$ws.Cells.Item(13, 6).Value = 1.028691572118302
$ws.Cells.Item(13, 9).Value = 1.029974951268907
$ws.Cells.Item(13, 10).Value = 1.025940748512562
$ws.Cells.Item(13, 11).Value = 1.028384137939082
$ws.Cells.Item(13, 12).Value = 1.026702599805097
$ws.Cells.Item(13, 13).Value = 1.032321605203291
$ws.Cells.Item(13, 14).Value = 1.027397702052498

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.019377694484398
$ws.Cells.Item(14, 4).Value = 1.02487676879935
$ws.Cells.Item(14, 5).Value = 1.023230404645906
$ws.Cells.Item(14, 6).Value = 1.028916111092129
$ws.Cells.Item(14, 9).Value = 1.030019183023224
$ws.Cells.Item(14, 10).Value = 1.026067518762001
$ws.Cells.Item(14, 11).Value = 1.02848937107966
$ws.Cells.Item(14, 12).Value = 1.026849276067554
$ws.Cells.Item(14, 13).Value = 1.032513426250582
$ws.Cells.Item(14, 14).Value = 1.027524652330232

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.019494220508949
$ws.Cells.Item(15, 4).Value = 1.024961641067027
$ws.Cells.Item(15, 5).Value = 1.023340964043608
$ws.Cells.Item(15, 6).Value = 1.029054449822677
$ws.Cells.Item(15, 9).Value = 1.030046388956613
$ws.Cells.Item(15, 10).Value = 1.026145608071956
$ws.Cells.Item(15, 11).Value = 1.028554180328548
$ws.Cells.Item(15, 12).Value = 1.026939637368016
$ws.Cells.Item(15, 13).Value = 1.032631593553246
$ws.Cells.Item(15, 14).Value = 1.027602852535964

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.020172364962748
$ws.Cells.Item(16, 4).Value = 1.02545544557324
$ws.Cells.Item(16, 5).Value = 1.023984524651965
$ws.Cells.Item(16, 6).Value = 1.029859533838605
$ws.Cells.Item(16, 9).Value = 1.030204026145488
$ws.Cells.Item(16, 10).Value = 1.026599846871533
$ws.Cells.Item(16, 11).Value = 1.028930965354474
$ws.Cells.Item(16, 12).Value = 1.027465415014644
$ws.Cells.Item(16, 13).Value = 1.033319074851329
$ws.Cells.Item(16, 14).Value = 1.028057736406739

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.020597671235268
$ws.Cells.Item(17, 4).Value = 1.025765029937762
$ws.Cells.Item(17, 5).Value = 1.024388264138395
$ws.Cells.Item(17, 6).Value = 1.03036444939224
$ws.Cells.Item(17, 9).Value = 1.030302279848624
$ws.Cells.Item(17, 10).Value = 1.026884537423329
$ws.Cells.Item(17, 11).Value = 1.029166931715504
$ws.Cells.Item(17, 12).Value = 1.027795077400509
$ws.Cells.Item(17, 13).Value = 1.03375004808083
$ws.Cells.Item(17, 14).Value = 1.028342831251773

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.020845718303057
$ws.Cells.Item(18, 4).Value = 1.025945545419294
$ws.Cells.Item(18, 5).Value = 1.024623777135707
$ws.Cells.Item(18, 6).Value = 1.030658925595058
$ws.Cells.Item(18, 9).Value = 1.030359363319266
$ws.Cells.Item(18, 10).Value = 1.027050505797744
$ws.Cells.Item(18, 11).Value = 1.029304429724301
$ws.Cells.Item(18, 12).Value = 1.027987312200266
$ws.Cells.Item(18, 13).Value = 1.034001332398519
$ws.Cells.Item(18, 14).Value = 1.028509035320317

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.020930291424572
$ws.Cells.Item(19, 4).Value = 1.026007086387231
$ws.Cells.Item(19, 5).Value = 1.024704084156193
$ws.Cells.Item(19, 6).Value = 1.030759328933356
$ws.Cells.Item(19, 9).Value = 1.030378788935878
$ws.Cells.Item(19, 10).Value = 1.027107082019723
$ws.Cells.Item(19, 11).Value = 1.029351289734219
$ws.Cells.Item(19, 12).Value = 1.028052850551801
$ws.Cells.Item(19, 13).Value = 1.034086997792251
$ws.Cells.Item(19, 14).Value = 1.028565691887019

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.020552042687856
$ws.Cells.Item(20, 4).Value = 1.025731820655662
$ws.Cells.Item(20, 5).Value = 1.024344944795149
$ws.Cells.Item(20, 6).Value = 1.030310280110666
$ws.Cells.Item(20, 9).Value = 1.030291761559462
$ws.Cells.Item(20, 10).Value = 1.026854001804818
$ws.Cells.Item(20, 11).Value = 1.029141628949243
$ws.Cells.Item(20, 12).Value = 1.027759713097918
$ws.Cells.Item(20, 13).Value = 1.033703818580463
$ws.Cells.Item(20, 14).Value = 1.028312252269182

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.019322000389209
$ws.Cells.Item(21, 4).Value = 1.024836201548532
$ws.Cells.Item(21, 5).Value = 1.023177564803843
$ws.Cells.Item(21, 6).Value = 1.028849991453131
$ws.Cells.Item(21, 9).Value = 1.030006167605161
$ws.Cells.Item(21, 10).Value = 1.026030191839336
$ws.Cells.Item(21, 11).Value = 1.028458388451734
$ws.Cells.Item(21, 12).Value = 1.026806085812232
$ws.Cells.Item(21, 13).Value = 1.032456943898198
$ws.Cells.Item(21, 14).Value = 1.027487272399056

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.018548372290175
$ws.Cells.Item(22, 4).Value = 1.024272550289347
$ws.Cells.Item(22, 5).Value = 1.02244374905419
$ws.Cells.Item(22, 6).Value = 1.027931540189022
$ws.Cells.Item(22, 9).Value = 1.029824563992348
$ws.Cells.Item(22, 10).Value = 1.025511443337651
$ws.Cells.Item(22, 11).Value = 1.02802756993953
$ws.Cells.Item(22, 12).Value = 1.026206032106899
$ws.Cells.Item(22, 13).Value = 1.031672113767756
$ws.Cells.Item(22, 14).Value = 1.026967787214996

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.01895851150088
$ws.Cells.Item(23, 4).Value = 1.024571403715004
$ws.Cells.Item(23, 5).Value = 1.022832743521904
$ws.Cells.Item(23, 6).Value = 1.028418458779641
$ws.Cells.Item(23, 9).Value = 1.029921028863848
$ws.Cells.Item(23, 10).Value = 1.025786516498964
$ws.Cells.Item(23, 11).Value = 1.028256072452904
$ws.Cells.Item(23, 12).Value = 1.026524176745019
$ws.Cells.Item(23, 13).Value = 1.032088250266598
$ws.Cells.Item(23, 14).Value = 1.02724325101175

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.020572660347022
$ws.Cells.Item(24, 4).Value = 1.02574682668549
$ws.Cells.Item(24, 5).Value = 1.024364518883983
$ws.Cells.Item(24, 6).Value = 1.030334756977052
$ws.Cells.Item(24, 9).Value = 1.030296515020538
$ws.Cells.Item(24, 10).Value = 1.026867799804398
$ws.Cells.Item(24, 11).Value = 1.02915306260383
$ws.Cells.Item(24, 12).Value = 1.027775692863468
$ws.Cells.Item(24, 13).Value = 1.033724707994936
$ws.Cells.Item(24, 14).Value = 1.028326069863503

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.022444740317305
$ws.Cells.Item(25, 4).Value = 1.027108492452672
$ws.Cells.Item(25, 5).Value = 1.026142790448508
$ws.Cells.Item(25, 6).Value = 1.032557248190971
$ws.Cells.Item(25, 9).Value = 1.030723370818177
$ws.Cells.Item(25, 10).Value = 1.028119165735809
$ws.Cells.Item(25, 11).Value = 1.030188589560907
$ws.Cells.Item(25, 12).Value = 1.029225982128866
$ws.Cells.Item(25, 13).Value = 1.035620003516613
$ws.Cells.Item(25, 14).Value = 1.02957921287807

Write-Output "Updated vm_pu values for 380 kV case (rows 2-25)"